# Applies the row-content changes described by the target diff.
# The worksheet rows themselves stay in place (row numbers unchanged); only the
# per-row data (observation id, taxon sort order, redlist status, species info,
# coordinates) is updated to match the target content for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 112044178
$ws.Range("B9").Value = 96735
$ws.Range("Q9").Value = 554737
$ws.Range("R9").Value = 6697621

# Row 10
$ws.Range("A10").Value = 112044172
$ws.Range("B10").Value = 90806
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 4361
$ws.Range("F10").Value = "Orange taggsvamp"
$ws.Range("G10").Value = "Hydnellum aurantiacum"
$ws.Range("H10").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q10").Value = 554722
$ws.Range("R10").Value = 6697604

# Row 11
$ws.Range("A11").Value = 112044176
$ws.Range("B11").Value = 96735
$ws.Range("Q11").Value = 554725
$ws.Range("R11").Value = 6697571

# Row 12
$ws.Range("A12").Value = 112044164
$ws.Range("B12").Value = 89072
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 256703
$ws.Range("F12").Value = "Tallfingersvamp"
$ws.Range("G12").Value = "Ramaria eosanguinea"
$ws.Range("H12").Value = "R.H.Petersen"
$ws.Range("Q12").Value = 554725
$ws.Range("R12").Value = 6697591

# Row 23
$ws.Range("A23").Value = 112044169
$ws.Range("B23").Value = 89993
$ws.Range("E23").Value = 1209
$ws.Range("F23").Value = "Rynkskinn"
$ws.Range("G23").Value = "Phlebia centrifuga"
$ws.Range("H23").Value = "P.Karst."
$ws.Range("Q23").Value = 554765
$ws.Range("R23").Value = 6697617

# Row 24
$ws.Range("A24").Value = 112044155
$ws.Range("B24").Value = 89553
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 1202
$ws.Range("F24").Value = "Ullticka"
$ws.Range("G24").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H24").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("R24").Value = 6697629

# Row 25
$ws.Range("A25").Value = 112044154
$ws.Range("B25").Value = 89553
$ws.Range("Q25").Value = 554768
$ws.Range("R25").Value = 6697637

# Row 26
$ws.Range("A26").Value = 112044157
$ws.Range("B26").Value = 89553
$ws.Range("D26").Value = "NT"
$ws.Range("E26").Value = 1202
$ws.Range("F26").Value = "Ullticka"
$ws.Range("G26").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H26").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q26").Value = 554764
$ws.Range("R26").Value = 6697616

# Row 27
$ws.Range("A27").Value = 112044156
$ws.Range("B27").Value = 89553
$ws.Range("Q27").Value = 554762
$ws.Range("R27").Value = 6697614

# Row 28
$ws.Range("B28").Value = 96735

# Row 29
$ws.Range("A29").Value = 112044168
$ws.Range("B29").Value = 89993
$ws.Range("D29").Value = "VU"
$ws.Range("E29").Value = 1209
$ws.Range("F29").Value = "Rynkskinn"
$ws.Range("G29").Value = "Phlebia centrifuga"
$ws.Range("H29").Value = "P.Karst."
$ws.Range("Q29").Value = 554761

# Row 30
$ws.Range("A30").Value = 112044179
$ws.Range("B30").Value = 96735
$ws.Range("D30").Value = "VU"
$ws.Range("E30").Value = 220787
$ws.Range("F30").Value = "Knärot"
$ws.Range("G30").Value = "Goodyera repens"
$ws.Range("H30").Value = "(L.) R. Br."
$ws.Range("Q30").Value = 554795
$ws.Range("R30").Value = 6697596

# AF10 had an empty placeholder cell before the edit; the diff removes it entirely.
$ws.Range("AF10").ClearContents()

# AF12 had no cell before the edit; the diff adds an empty placeholder cell there.
$ws.Range("AF12").NumberFormat = "@"

Write-Output "Edit complete"
